$p = $ppt.ActivePresentation

# -----------------------------------------------------------------
# Slide 5: "Para os seminários, escolher 5 assuntos entre:" ->
#          "Para os seminários, escolher 5 assuntos entre (sugestão):"
# The final text is split into 3 runs:
#   "Para os seminários, escolher 5 assuntos entre"
#   " (sugestão)"
#   ":"
# -----------------------------------------------------------------
$s5 = $p.Slides.Item(5)
$shSem = $s5.Shapes.Item(5)
$trSem = $shSem.TextFrame.TextRange
$paraSem = $trSem.Paragraphs(7, 1)

# Replace the trailing colon with " (sugestão):" -- this splits the
# original single run into "...entre" + " (sugestão):"
$colon = $paraSem.Characters($paraSem.Text.TrimEnd().Length, 1)
$colon.Text = " (sugestão):"

# Re-fetch the paragraph (its length changed) and split the trailing
# ":" back out into its own run by re-assigning its own text.
$paraSem = $trSem.Paragraphs(7, 1)
$lastLen = $paraSem.Text.TrimEnd().Length
$tail = $paraSem.Characters($lastLen, 1)
$tail.Text = ":"

# -----------------------------------------------------------------
# Slide 6: update the class schedule times from :30 to :00
#   "8:30h – 12:30h."   -> "8:00h – 12:00h."
#   "13:30 – 17:30h"    -> "13:00 – 17:00h"
# Each paragraph ends up split into 5 runs around the two "00"
# replacements (matching how PowerPoint splits runs on edit).
# -----------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shHora = $s6.Shapes.Item(4)
$trHora = $shHora.TextFrame.TextRange

# Paragraph 2: "8:30h – 12:30h."
$paraManha = $trHora.Paragraphs(2, 1)
$txt = $paraManha.Text
$firstIdx = $txt.IndexOf("30") + 1
$sub1 = $paraManha.Characters($firstIdx, 2)
$sub1.Text = "00"

$paraManha = $trHora.Paragraphs(2, 1)
$txt = $paraManha.Text
$secondIdx = $txt.IndexOf("30") + 1
$sub2 = $paraManha.Characters($secondIdx, 2)
$sub2.Text = "00"

# Paragraph 4: "13:30 – 17:30h"
$paraTarde = $trHora.Paragraphs(4, 1)
$txt2 = $paraTarde.Text
$firstIdx2 = $txt2.IndexOf("30") + 1
$sub3 = $paraTarde.Characters($firstIdx2, 2)
$sub3.Text = "00"

$paraTarde = $trHora.Paragraphs(4, 1)
$txt2 = $paraTarde.Text
$secondIdx2 = $txt2.IndexOf("30") + 1
$sub4 = $paraTarde.Characters($secondIdx2, 2)
$sub4.Text = "00"
